$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Hedieh Eshaghi (row 5) grade columns F:K
$ws.Range("F5").Value = 100
$ws.Range("G5").Value = 100
$ws.Range("H5").Value = 100
$ws.Range("I5").Value = 95
$ws.Range("J5").Value = 100
$ws.Range("K5").Value = 100

# Update the active selection on the sheet view
$ws.Range("J6").Select()
